# Implement hh-based ACF around TB cases detected through PCF
#
# - bump the "population" parameter (B1) from 1000 to 10000
# - append a new parameter row "hh_based_acf_coverage_perc" at the
#   bottom of the `constant` sheet (row 79), formatted like the other
#   red/yellow-highlighted parameter rows (e.g. row 43), with a 0
#   default value and type "float"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constant")

# population: 1000 -> 10000
$ws.Range("B1").Value = 10000

# Grab the formatting of an existing "highlighted parameter" row (A43:C43 -
# perc_smearpos) and stamp it onto the new row 79 so the new parameter
# matches the house style used throughout the sheet.
$ws.Range("A43:C43").Copy()
$ws.Range("A79:C79").PasteSpecial(-4122)

$ws.Range("A79").Value = "hh_based_acf_coverage_perc"
$ws.Range("B79").Value = 0
$ws.Range("C79").Value = "float"

# Match the new selection/cursor position left by the edit.
$ws.Range("B79").Select() | Out-Null
